$wb = $excel.ActiveWorkbook

# --- Fix selection on the Allegany County sheet (no longer the last-active tab) ---
$wsAllegany = $wb.Worksheets.Item("Allegany County")
$wsAllegany.Activate()
$wsAllegany.Range("A1:D1").Select() | Out-Null

# --- Add the new Montgomery County sheet at the end of the tab strip ---
$count = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($count))
$ws.Name = "Montgomery County"

# Header row
$ws.Cells.Item(1,1).Value = 'Zone'
$ws.Cells.Item(1,2).Value = 'Zone Abbreviation'
$ws.Cells.Item(1,3).Value = 'Issuing Body'
$ws.Cells.Item(1,4).Value = 'Zone General Description'

# Data rows (Montgomery County zoning table)
$ws.Cells.Item(2,1).Value = 'Agricultural Reserve'
$ws.Cells.Item(2,2).Value = 'AR'
$ws.Cells.Item(2,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(2,4).Value = 'The intent of the AR zone is to promote agriculture as the primary land use in areas of the County designated for agricultural preservation in the general plan, the Functional Master Plan for Preservation of Agriculture and Rural Open Space, and other current or future master plans. The AR zone accomplishes this intent by providing large areas of generally contiguous properties suitable for agricultural and related uses and permitting the transfer of development rights from properties in this zone to properties in designated receiving areas.'
$ws.Cells.Item(3,1).Value = 'Rural'
$ws.Cells.Item(3,2).Value = 'R'
$ws.Cells.Item(3,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(3,4).Value = 'The intent of the R zone is to preserve rural areas of the County for agriculture and other natural resource development,residential uses of a rural character, extensive recreational facilities, and protection of scenic and environmentally sensitive areas.'
$ws.Cells.Item(4,1).Value = 'Rural Cluster'
$ws.Cells.Item(4,2).Value = 'RC'
$ws.Cells.Item(4,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(4,4).Value = 'The intent of the RC zone is to provide designated areas of the County for a compatible mixture of agricultural uses and very low-density residential development, to promote agriculture, and to protect scenic and environmentally sensitive areas. The RC zone permits an optional method Cluster Development alternative to provide greater flexibility in achieving a compatible mixture of agricultural and residential uses and to protect scenic and environmentally sensitive areas without jeopardizing farming or other agricultural uses.'
$ws.Cells.Item(5,1).Value = 'Rural Neighborhood Cluster'
$ws.Cells.Item(5,2).Value = 'RNC'
$ws.Cells.Item(5,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(5,4).Value = 'The intent of the RNC zone is to preserve open land, environmentally sensitive natural resources, and rural community character through clustering of residential development in the form of small neighborhoods that provide neighborhood identity in an open space setting. A master plan must recommend the RNC zone, and must provide development guidelines and recommendations regarding density, and the location and rationale for preserving the rural open space. It is also the intent of the RNC zone to implement the recommendations of the applicable master plan, such as maintaining broad vistas of open space, preserving agrarian character, or preserving environmentally sensitive natural resources to the maximum extent possible, and to ensure that new development is in harmony with the policies and guidelines of the applicable master plan and is compatible with existing development in adjoining communities.'
$ws.Cells.Item(6,1).Value = 'Residential Estate - 2'
$ws.Cells.Item(6,2).Value = 'RE-2'
$ws.Cells.Item(6,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(6,4).Value = 'The intent of the RE-2 zone is to provide designated areas of the County for large-lot residential uses. The predominant use is residential in a detached house.'
$ws.Cells.Item(7,1).Value = 'Residential Estate - 2C'
$ws.Cells.Item(7,2).Value = 'RE-2C'
$ws.Cells.Item(7,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(7,4).Value = 'The intent of the RE-2C zone is to provide designated areas of the County for large-lot residential uses. The predominant use is residential in a detached house.'
$ws.Cells.Item(8,1).Value = 'Residential Estate - 1'
$ws.Cells.Item(8,2).Value = 'RE-1'
$ws.Cells.Item(8,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(8,4).Value = 'The intent of the RE-1 zone is to provide designated areas of the County for large-lot residential uses. The predominant use is residential in a detached house.'
$ws.Cells.Item(9,1).Value = 'Residential - 200'
$ws.Cells.Item(9,2).Value = 'R-200'
$ws.Cells.Item(9,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(9,4).Value = 'The intent of the R-200 zone is to provide designated areas of the County for residential uses with a minimum lot size of 20,000 square feet. The predominant use is residential in a detached house.'
$ws.Cells.Item(10,1).Value = 'Residential - 90'
$ws.Cells.Item(10,2).Value = 'R-90'
$ws.Cells.Item(10,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(10,4).Value = 'The intent of the R-90 zone is to provide designated areas of the County for moderate density residential uses. The predominant use is residential in a detached house. A limited number of other building types may be allowed under the optional method of development'
$ws.Cells.Item(11,1).Value = 'Residential - 60'
$ws.Cells.Item(11,2).Value = 'R-60'
$ws.Cells.Item(11,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(11,4).Value = 'The intent of the R-60 zone is to provide designated areas of the County for moderate density residential uses. The predominant use is residential in a detached house. A limited number of other building types may be allowed under the optional method of development.'
$ws.Cells.Item(12,1).Value = 'Residential - 40'
$ws.Cells.Item(12,2).Value = 'R- 40'
$ws.Cells.Item(12,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(12,4).Value = 'The intent of the R-40 zone is to provide designated areas of the County for moderate density residential uses. The predominant use is residential in a duplex or detached house. A limited number of other building types may be allowed under the optional method of development.'
$ws.Cells.Item(13,1).Value = 'Townhouse Low Density'
$ws.Cells.Item(13,2).Value = 'TLD'
$ws.Cells.Item(13,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(13,4).Value = 'The intent of the TLD zone is to provide designated areas of the County for residential purposes at slightly higher densities than the R-90, R-60, and R-40 zones. It is also the intent of the TLD zone to provide a buffer or transition between nonresidential or high-density residential uses and the medium- or low-density Residential zones.'
$ws.Cells.Item(14,1).Value = 'Townhouse Medium Density'
$ws.Cells.Item(14,2).Value = 'TMD'
$ws.Cells.Item(14,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(14,4).Value = 'The intent of the TMD zone is to provide designated areas of the County for residential purposes at slightly higher densities than the R-90, R-60, and R-40 zones. It is also the intent of the TMD zone to provide a buffer or transition between nonresidential or high-density residential uses and the medium- or low-density Residential zones.'
$ws.Cells.Item(15,1).Value = 'Townhouse High Density'
$ws.Cells.Item(15,2).Value = 'THD'
$ws.Cells.Item(15,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(15,4).Value = 'The intent of the THD zone is to provide designated areas of the County for residential purposes at slightly higher densities than the R-90, R-60, and R-40 zones. It is also the intent of the THD zone to provide a buffer or transition between nonresidential or high-density residential uses and the medium- or low-density Residential zones.'
$ws.Cells.Item(16,1).Value = 'Residentual Multi-Unit Low Density'
$ws.Cells.Item(16,2).Value = 'R-30'
$ws.Cells.Item(16,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(16,4).Value = 'The intent of the R-30 zone is to provide designated areas of the County for higher-density, multi-unit residential uses. The predominant use is residential in an apartment building, although detached house, duplex, and townhouse building types are allowed.'
$ws.Cells.Item(17,1).Value = 'Residential Multi-Unit Medium Density'
$ws.Cells.Item(17,2).Value = 'R-20'
$ws.Cells.Item(17,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(17,4).Value = 'The intent of the R-20 zone is to provide designated areas of the County for higher-density, multi-unit residential uses. The predominant use is residential in an apartment building, although detached house, duplex, and townhouse building types are allowed.'
$ws.Cells.Item(18,1).Value = 'Residential Multi-Unit High Density'
$ws.Cells.Item(18,2).Value = 'R-10'
$ws.Cells.Item(18,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(18,4).Value = 'The intent of the R-10 zone is to provide designated areas of the County for higher-density, multi-unit residential uses. The predominant use is residential in an apartment building, although detached house, duplex, and townhouse building types are allowed.'
$ws.Cells.Item(19,1).Value = 'Commercial Residential Neighborhood'
$ws.Cells.Item(19,2).Value = 'CRN'
$ws.Cells.Item(19,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(19,4).Value = 'The CRN zone is intended for pedestrian-scale, neighborhood-serving mixed-use centers and transitional edges. Retail tenant ground floor footprints are limited to preserve community scale.'
$ws.Cells.Item(20,1).Value = 'Commercial Residential Town'
$ws.Cells.Item(20,2).Value = 'CRT'
$ws.Cells.Item(20,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(20,4).Value = 'The CRT zone is intended for small downtown, mixed-use, pedestrian-oriented centers and edges of larger, more intense downtowns. Retail tenant ground floor footprints are limited to preserve the town center scale. Transit options may include light rail, Metro, and bus.'
$ws.Cells.Item(21,1).Value = 'Commercial Residential'
$ws.Cells.Item(21,2).Value = 'CR'
$ws.Cells.Item(21,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(21,4).Value = 'The CR zone is intended for larger downtown, mixed-use, and pedestrian-oriented areas in close proximity to transit options such as Metro, light rail, and bus. Retail tenant gross floor area is not restricted.'
$ws.Cells.Item(22,1).Value = 'General Retail'
$ws.Cells.Item(22,2).Value = 'GR'
$ws.Cells.Item(22,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(22,4).Value = 'The GR zone is intended for commercial areas of a general nature, including regional shopping centers and clusters of commercial development. The GR zone provides development opportunities adjacent to the County''s most auto-dominated corridors and those areas with few alternative mobility options. The GR zone allows flexibility in building, circulation, and parking lot layout. Retail/Service Establishment gross floor area is not restricted.'
$ws.Cells.Item(23,1).Value = 'Neighborhood Retail'
$ws.Cells.Item(23,2).Value = 'NR'
$ws.Cells.Item(23,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(23,4).Value = 'The NR zone is intended for commercial areas that have a neighborhood orientation and which supply necessities usually requiring frequent purchasing and convenient automobile access. The NR zone addresses development opportunities within primarily residential areas with few alternative mobility options and without a critical mass of density needed for pedestrian-oriented commercial uses. The NR zone allows flexibility in building, circulation, and parking lot layout.'
$ws.Cells.Item(24,1).Value = 'Life Sciences Center'
$ws.Cells.Item(24,2).Value = 'LSC'
$ws.Cells.Item(24,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(24,4).Value = 'The LSC zone is intended primarily for research, development, education, and related activities. The primary purpose is to promote research, academic, and clinical facilities that advance the life sciences, health care services, and applied technologies. It is also the purpose of the LSC zone to provide opportunities for the development of uses that support a Life Sciences Center while retaining an environment conducive to high technology research, development, and production. Retail sales and personal services are allowed but are intended for the convenience of employees and residents in the zone.'
$ws.Cells.Item(25,1).Value = 'Employment Office'
$ws.Cells.Item(25,2).Value = 'EOF'
$ws.Cells.Item(25,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(25,4).Value = 'The EOF zone is intended for office and employment activity combined with limited residential and neighborhood commercial uses. The EOF allows flexibility in building, circulation, and parking lot layout.'
$ws.Cells.Item(26,1).Value = 'Light Industrial'
$ws.Cells.Item(26,2).Value = 'IL'
$ws.Cells.Item(26,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(26,4).Value = 'The IL zone is intended to provide land for industrial activities where major transportation links are not typically necessary and noise, dust, vibration, glare, odors, and other adverse environmental impacts are usually minimal. The IL zone is appropriate as a transitional Industrial zone between a Residentially zoned area and land classified in the IM and IH zones.'
$ws.Cells.Item(27,1).Value = 'Moderate Industrial'
$ws.Cells.Item(27,2).Value = 'IM'
$ws.Cells.Item(27,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(27,4).Value = 'The IM zone is intended to provide land for industrial activities where major transportation links are not typically necessary and noise, dust, vibration, glare, odors, and other adverse environmental impacts are usually minimal.'
$ws.Cells.Item(28,1).Value = 'Heavy Industrial'
$ws.Cells.Item(28,2).Value = 'IH'
$ws.Cells.Item(28,3).Value = 'Montgomery County Planning and Zoning Department'
$ws.Cells.Item(28,4).Value = 'The IH zone is intended to provide land for industrial activities that usually need major transportation links to highways or rail and may create significant noise, dust, vibration, glare, odors, and other adverse environmental impacts.'

# Autofit columns A:D based on content
$ws.Range("A1:D28").EntireColumn.AutoFit() | Out-Null

# Apply Corbel/#212529 font styling to the Zone General Description column (D),
# matching the style used for this column throughout the workbook.
$ws.Cells.Item(2,4).Font.Name = "Corbel"
$ws.Cells.Item(2,4).Font.Size = 11
$ws.Cells.Item(2,4).Font.Color = 2696481
$ws.Cells.Item(2,4).Copy()
$ws.Range("D3:D28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Select D28 (last cell of data) and make this the active sheet/tab,
# matching the final view state after entering the data.
$ws.Activate()
$ws.Cells.Item(28,4).Select() | Out-Null
